# Actualizacion automatica del tracker
# - Completa el resultado pendiente de la fila 19 (Acierto / +2.5)
# - Agrega 5 filas nuevas (26-30) con los partidos del 2025-09-20

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Resolver resultado pendiente de la fila 19 ---
$ws.Range("G19").Value = "Acierto"
$ws.Range("H19").Value = 2.5

# --- 2) Helper para escribir una fila nueva del tracker ---
# (usa parametros posicionales: el enlace de parametros con nombre
#  "-Nombre valor" no funciona de forma fiable en este runtime)
function Add-TrackerRow {
    param($Row, $EventId, $Fecha, $JugadorA, $JugadorB, $Pronostico, $Cuota)

    $ws.Range("A$Row").Value = $EventId

    # La columna "fecha" guarda texto plano (p.ej. "2025-09-20"), no una
    # fecha real, asi que forzamos formato de texto antes de escribir el
    # valor para que Excel no lo autoconvierta a un numero de fecha, y
    # luego restauramos el estilo por defecto de la celda.
    $ws.Range("B$Row").NumberFormat = "@"
    $ws.Range("B$Row").Value = $Fecha
    $ws.Range("B$Row").Style = "Normal"

    $ws.Range("C$Row").Value = $JugadorA
    $ws.Range("D$Row").Value = $JugadorB
    $ws.Range("E$Row").Value = $Pronostico
    $ws.Range("F$Row").Value = $Cuota

    # "resultado" y "profit" quedan pendientes (cadena vacia) hasta que se
    # conozca el desenlace del partido. Se usa la comilla simple (prefijo
    # de texto) para forzar una cadena vacia en vez de una celda en blanco,
    # y luego se restaura el estilo por defecto.
    $ws.Range("G$Row").Value = "'"
    $ws.Range("G$Row").Style = "Normal"
    $ws.Range("H$Row").Value = "'"
    $ws.Range("H$Row").Style = "Normal"
}

Add-TrackerRow 26 14722601 "2025-09-20" "Jiaqi Wang" "Ying Zhang" "Gana Jiaqi Wang" 2.5
Add-TrackerRow 27 14722597 "2025-09-20" "Riya Bhatia" "Ankita Raina" "Gana Riya Bhatia" 2.1
Add-TrackerRow 28 14722598 "2025-09-20" "Diletta Cherubini" "Jing-Jing Lu" "Gana Jing-Jing Lu" 2.2
Add-TrackerRow 29 14722603 "2025-09-20" "Ye Qiuyu" "Hong Yi Cody Wong" "Gana Ye Qiuyu" 2.38
Add-TrackerRow 30 14722602 "2025-09-20" "Lea Ma" "Fang Ying Xun" "Gana Fang Ying Xun" 2.5

Write-Host "Tracker actualizado: fila 19 resuelta y filas 26-30 agregadas."
